# Commit: "add price data and check"
#
# 1. Rename the price-reference sheet "p_ref" -> "price_ref".
# 2. Add the price-derived cost data: cost_ref!B2 100 -> 15, with an
#    updated note explaining the value was reduced from COST_NODAL_NET
#    (was previously a placeholder note saying COST_NODAL_NET was empty).
# 3. Leave the cost_ref sheet's selection on the note cell (C2).
# 4. Make "price_ref" the active sheet/tab with G20 selected (previously
#    "gdp_calibrate" was the active tab).

$wb = $excel.ActiveWorkbook

# --- rename the price-reference sheet ---
$priceSheet = $wb.Worksheets.Item("p_ref")
$priceSheet.Name = "price_ref"

# --- add the price data / check on cost_ref ---
$costSheet = $wb.Worksheets.Item("cost_ref")
$costSheet.Range("B2").Value = 15
$costSheet.Range("C2").Value = "reduced from value found in COST_NODAL_NET"
$costSheet.Range("C2").Select()

# --- price_ref becomes the active tab / selection ---
$priceSheet.Activate()
$priceSheet.Range("G20").Select()
